# Auto-generated edit script
# Applies updated market price data to the Chocobo profit-tracking workbook
# (H/I/J/K/L/M/N columns: price + profit calculations per crafting leve)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5887187
$ws.Range("I74").Value = 8337683.5
$ws.Range("K74").Value = 8337683.5
$ws.Range("M74").Value = -8336747.5
$ws.Range("H77").Value = 5887187
$ws.Range("I77").Value = 8337683.5
$ws.Range("K77").Value = 41688417.5
$ws.Range("M77").Value = -41683737.5
$ws.Range("H98").Value = 6491.2856
$ws.Range("I98").Value = 3935.75
$ws.Range("J98").Value = 9898.666999999999
$ws.Range("K98").Value = 3935.75
$ws.Range("L98").Value = 9898.666999999999
$ws.Range("M98").Value = -2437.75
$ws.Range("N98").Value = -12894.667
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882
$ws.Range("H122").Value = 6491.2856
$ws.Range("I122").Value = 3935.75
$ws.Range("J122").Value = 9898.666999999999
$ws.Range("K122").Value = 11807.25
$ws.Range("L122").Value = 29696.001
$ws.Range("M122").Value = -9357.25
$ws.Range("N122").Value = -34596.001
$ws.Range("H131").Value = 6400
$ws.Range("I131").Value = 4600
$ws.Range("J131").Value = 7750
$ws.Range("K131").Value = 13800
$ws.Range("L131").Value = 23250
$ws.Range("M131").Value = -8760
$ws.Range("N131").Value = -33330
$ws.Range("H134").Value = 40891.31
$ws.Range("J134").Value = 40891.31
$ws.Range("L134").Value = 40891.31
$ws.Range("N134").Value = -51031.31
$ws.Range("H137").Value = 4600.5625
$ws.Range("I137").Value = 4799.35
$ws.Range("J137").Value = 4269.25
$ws.Range("K137").Value = 14398.05
$ws.Range("L137").Value = 12807.75
$ws.Range("M137").Value = -11848.05
$ws.Range("N137").Value = -17907.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1342.4166
$ws.Range("I61").Value = 1085.6666
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 1085.6666
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -873.6666
$ws.Range("N61").Value = -4590.6665
$ws.Range("H74").Value = 3302.6667
$ws.Range("I74").Value = 3191.457
$ws.Range("J74").Value = 3858.7144
$ws.Range("K74").Value = 3191.457
$ws.Range("L74").Value = 3858.7144
$ws.Range("M74").Value = -2317.457
$ws.Range("N74").Value = -5606.7144
$ws.Range("H77").Value = 3302.6667
$ws.Range("I77").Value = 3191.457
$ws.Range("J77").Value = 3858.7144
$ws.Range("K77").Value = 15957.285
$ws.Range("L77").Value = 19293.572
$ws.Range("M77").Value = -11589.285
$ws.Range("N77").Value = -28029.572
$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344
$ws.Range("H110").Value = 992.6667
$ws.Range("I110").Value = 936.36365
$ws.Range("J110").Value = 1147.5
$ws.Range("K110").Value = 936.36365
$ws.Range("L110").Value = 1147.5
$ws.Range("M110").Value = 1108.63635
$ws.Range("N110").Value = -5237.5
$ws.Range("H136").Value = 1342.4166
$ws.Range("I136").Value = 1085.6666
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 3256.9998
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -706.9998000000001
$ws.Range("N136").Value = -17599.9995
$ws.Range("H137").Value = 40595
$ws.Range("J137").Value = 40595
$ws.Range("L137").Value = 40595
$ws.Range("N137").Value = -50795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 38960
$ws.Range("J137").Value = 40511.11
$ws.Range("L137").Value = 40511.11
$ws.Range("N137").Value = -50711.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10206649
$ws.Range("I31").Value = 1299.4412
$ws.Range("J31").Value = 33338774
$ws.Range("K31").Value = 1299.4412
$ws.Range("L31").Value = 33338774
$ws.Range("M31").Value = -1004.4412
$ws.Range("N31").Value = -33339364
$ws.Range("H34").Value = 10206649
$ws.Range("I34").Value = 1299.4412
$ws.Range("J34").Value = 33338774
$ws.Range("K34").Value = 1299.4412
$ws.Range("L34").Value = 33338774
$ws.Range("M34").Value = -1097.4412
$ws.Range("N34").Value = -33339178
$ws.Range("H39").Value = 19595.182
$ws.Range("I39").Value = 4666.6665
$ws.Range("J39").Value = 25193.375
$ws.Range("K39").Value = 4666.6665
$ws.Range("L39").Value = 25193.375
$ws.Range("M39").Value = -4275.6665
$ws.Range("N39").Value = -25975.375
$ws.Range("H49").Value = 19595.182
$ws.Range("I49").Value = 4666.6665
$ws.Range("J49").Value = 25193.375
$ws.Range("K49").Value = 4666.6665
$ws.Range("L49").Value = 25193.375
$ws.Range("M49").Value = -4484.6665
$ws.Range("N49").Value = -25557.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 26.9375
$ws.Range("I12").Value = 6.5
$ws.Range("J12").Value = 33.75
$ws.Range("K12").Value = 19.5
$ws.Range("L12").Value = 101.25
$ws.Range("M12").Value = 153.5
$ws.Range("N12").Value = -447.25
$ws.Range("H34").Value = 11304.35
$ws.Range("J34").Value = 7233.5625
$ws.Range("L34").Value = 21700.6875
$ws.Range("N34").Value = -21868.6875
$ws.Range("H39").Value = 9971.619000000001
$ws.Range("J39").Value = 9971.619000000001
$ws.Range("L39").Value = 29914.857
$ws.Range("N39").Value = -30502.857
$ws.Range("H55").Value = 4752.4
$ws.Range("J55").Value = 4867.0835
$ws.Range("L55").Value = 14601.2505
$ws.Range("N55").Value = -14955.2505
$ws.Range("H75").Value = 2902.6
$ws.Range("I75").Value = 1013
$ws.Range("J75").Value = 3375
$ws.Range("K75").Value = 3039
$ws.Range("L75").Value = 10125
$ws.Range("M75").Value = -2041
$ws.Range("N75").Value = -12121
$ws.Range("H78").Value = 2902.6
$ws.Range("I78").Value = 1013
$ws.Range("J78").Value = 3375
$ws.Range("K78").Value = 9117
$ws.Range("L78").Value = 30375
$ws.Range("M78").Value = -4125
$ws.Range("N78").Value = -40359
$ws.Range("H103").Value = 1750
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 2166.6667
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 6500.000100000001
$ws.Range("M103").Value = -621
$ws.Range("N103").Value = -8258.000100000001
$ws.Range("H107").Value = 560.8261
$ws.Range("I107").Value = 406.46155
$ws.Range("J107").Value = 761.5
$ws.Range("K107").Value = 1219.38465
$ws.Range("L107").Value = 2284.5
$ws.Range("M107").Value = 700.61535
$ws.Range("N107").Value = -6124.5
$ws.Range("H131").Value = 6173711
$ws.Range("I131").Value = 83335220
$ws.Range("J131").Value = 790.5467
$ws.Range("K131").Value = 250005660
$ws.Range("L131").Value = 2371.6401
$ws.Range("M131").Value = -250000620
$ws.Range("N131").Value = -12451.6401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5258.5
$ws.Range("I36").Value = 1017
$ws.Range("J36").Value = 9500
$ws.Range("K36").Value = 1017
$ws.Range("L36").Value = 9500
$ws.Range("M36").Value = -532
$ws.Range("N36").Value = -10470
$ws.Range("H137").Value = 50169.285
$ws.Range("J137").Value = 57670
$ws.Range("L137").Value = 57670
$ws.Range("N137").Value = -67870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 26500
$ws.Range("J112").Value = 26500
$ws.Range("L112").Value = 26500
$ws.Range("N112").Value = -29454
$ws.Range("H122").Value = 5591.6665
$ws.Range("I122").Value = 3316.6667
$ws.Range("J122").Value = 7866.6665
$ws.Range("K122").Value = 9950.000100000001
$ws.Range("L122").Value = 23599.9995
$ws.Range("M122").Value = -7500.000100000001
$ws.Range("N122").Value = -28499.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8131732.5
$ws.Range("I132").Value = 1023.31036
$ws.Range("J132").Value = 27780946
$ws.Range("K132").Value = 3069.93108
$ws.Range("L132").Value = 83342838
$ws.Range("M132").Value = -539.9310799999998
$ws.Range("N132").Value = -83347898
$ws.Range("H135").Value = 41715
$ws.Range("J135").Value = 41715
$ws.Range("L135").Value = 41715
$ws.Range("N135").Value = -51855
